$wb = $excel.ActiveWorkbook

# --- Sheet: Means ---
$wsMeans = $wb.Worksheets.Item("Means")

# Row 9: Total Cancer Risk (per million)
$wsMeans.Range("B9").Value = 26
$wsMeans.Range("C9").Value = 26
$wsMeans.Range("D9").Value = 28
$wsMeans.Range("E9").Value = 29
$wsMeans.Range("F9").Value = 29
$wsMeans.Range("G9").Value = 29

# Row 10: Total Respiratory (hazard quotient)
$wsMeans.Range("B10").Value = 0.31
$wsMeans.Range("C10").Value = 0.32
$wsMeans.Range("D10").Value = 0.33
$wsMeans.Range("E10").Value = 0.36
$wsMeans.Range("F10").Value = 0.37
$wsMeans.Range("G10").Value = 0.35

# --- Sheet: Standard Deviations ---
$wsSD = $wb.Worksheets.Item("Standard Deviations")

# Row 9: Total Cancer Risk (per million) SD
$wsSD.Range("B9").Value = 8.3
$wsSD.Range("C9").Value = 5.4
$wsSD.Range("D9").Value = 4.1
$wsSD.Range("E9").Value = 3.1
$wsSD.Range("F9").Value = 2.3
$wsSD.Range("G9").Value = 2.6

# Row 10: Total Respiratory (hazard quotient) SD
$wsSD.Range("B10").Value = 0.11
$wsSD.Range("C10").Value = 0.043
$wsSD.Range("D10").Value = 0.041
$wsSD.Range("E10").Value = 0.049
$wsSD.Range("F10").Value = 0.042
$wsSD.Range("G10").Value = 0.048
